$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("S2").Value = 1.73
$ws.Range("T2").Value = 2.08

# Row 4 updates
$ws.Range("Q4").Value = 2.2
$ws.Range("R4").Value = 1.65
